$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 59.45197733333333
$ws.Range("H2").Value = 178.355932
$ws.Range("I2").Value = 0.304222453049858
$ws.Range("J2").Value = 0.304222453049858
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 94.95332266666666
$ws.Range("N2").Value = 284.859968
$ws.Range("O2").Value = 0.6912729447872573
$ws.Range("P2").Value = 0.6912729447872574
$ws.Range("Q2").Value = 5645.162786903353
$ws.Range("R2").Value = 50806.46508213017
$ws.Range("S2").Value = 0.2103007509901784
$ws.Range("T2").Value = 0.2103007509901785
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 59.45197733333333
$ws.Range("H3").Value = 178.355932
$ws.Range("I3").Value = 0.304222453049858
$ws.Range("J3").Value = 0.304222453049858
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 36.59611566666667
$ws.Range("N3").Value = 109.788347
$ws.Range("O3").Value = 0.2664246382770613
$ws.Range("P3").Value = 0.2664246382770613
$ws.Range("Q3").Value = 2175.711439102712
$ws.Range("R3").Value = 19581.4029519244
$ws.Range("S3").Value = 0.08105235700956866
$ws.Range("T3").Value = 0.08105235700956867
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 59.45197733333333
$ws.Range("H4").Value = 178.355932
$ws.Range("I4").Value = 0.304222453049858
$ws.Range("J4").Value = 0.304222453049858
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.810664333333333
$ws.Range("N4").Value = 17.431993
$ws.Range("O4").Value = 0.04230241693568138
$ws.Range("P4").Value = 0.04230241693568138
$ws.Range("Q4").Value = 345.4554842369417
$ws.Range("R4").Value = 3109.099358132476
$ws.Range("S4").Value = 0.01286934505011084
$ws.Range("T4").Value = 0.01286934505011085
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 107.1770123333333
$ws.Range("H5").Value = 321.531037
$ws.Range("I5").Value = 0.548436823552382
$ws.Range("J5").Value = 0.5484368235523819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 94.95332266666666
$ws.Range("N5").Value = 284.859968
$ws.Range("O5").Value = 0.6912729447872573
$ws.Range("P5").Value = 0.6912729447872574
$ws.Range("Q5").Value = 10176.81343453631
$ws.Range("R5").Value = 91591.3209108268
$ws.Range("S5").Value = 0.3791195380468245
$ws.Range("T5").Value = 0.3791195380468245
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 107.1770123333333
$ws.Range("H6").Value = 321.531037
$ws.Range("I6").Value = 0.548436823552382
$ws.Range("J6").Value = 0.5484368235523819
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 36.59611566666667
$ws.Range("N6").Value = 109.788347
$ws.Range("O6").Value = 0.2664246382770613
$ws.Range("P6").Value = 0.2664246382770613
$ws.Range("Q6").Value = 3922.262340158426
$ws.Range("R6").Value = 35300.36106142584
$ws.Range("S6").Value = 0.1461170823327639
$ws.Range("T6").Value = 0.1461170823327639
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 107.1770123333333
$ws.Range("H7").Value = 321.531037
$ws.Range("I7").Value = 0.548436823552382
$ws.Range("J7").Value = 0.5484368235523819
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.810664333333333
$ws.Range("N7").Value = 17.431993
$ws.Range("O7").Value = 0.04230241693568138
$ws.Range("P7").Value = 0.04230241693568138
$ws.Range("Q7").Value = 622.7696429185266
$ws.Range("R7").Value = 5604.92678626674
$ws.Range("S7").Value = 0.02320020317279358
$ws.Range("T7").Value = 0.02320020317279358
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.793724
$ws.Range("H8").Value = 86.38117199999999
$ws.Range("I8").Value = 0.1473407233977601
$ws.Range("J8").Value = 0.1473407233977601
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 94.95332266666666
$ws.Range("N8").Value = 284.859968
$ws.Range("O8").Value = 0.6912729447872573
$ws.Range("P8").Value = 0.6912729447872574
$ws.Range("Q8").Value = 2734.059765746944
$ws.Range("R8").Value = 24606.53789172249
$ws.Range("S8").Value = 0.1018526557502544
$ws.Range("T8").Value = 0.1018526557502544
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.793724
$ws.Range("H9").Value = 86.38117199999999
$ws.Range("I9").Value = 0.1473407233977601
$ws.Range("J9").Value = 0.1473407233977601
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 36.59611566666667
$ws.Range("N9").Value = 109.788347
$ws.Range("O9").Value = 0.2664246382770613
$ws.Range("P9").Value = 0.2664246382770613
$ws.Range("Q9").Value = 1053.738453978076
$ws.Range("R9").Value = 9483.646085802684
$ws.Range("S9").Value = 0.03925519893472876
$ws.Range("T9").Value = 0.03925519893472876
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.793724
$ws.Range("H10").Value = 86.38117199999999
$ws.Range("I10").Value = 0.1473407233977601
$ws.Range("J10").Value = 0.1473407233977601
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.810664333333333
$ws.Range("N10").Value = 17.431993
$ws.Range("O10").Value = 0.04230241693568138
$ws.Range("P10").Value = 0.04230241693568138
$ws.Range("Q10").Value = 167.3106650706439
$ws.Range("R10").Value = 1505.795985635796
$ws.Range("S10").Value = 0.006232868712776951
$ws.Range("T10").Value = 0.006232868712776951
